$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D6").Value = -7.959000000000001
$ws.Range("B7").Value = 5.815
$ws.Range("E7").Value = 16.604
$ws.Range("A9").Value = -21.752
$ws.Range("E10").Value = 16.492
$ws.Range("B12").Value = 5.927000000000001
$ws.Range("E13").Value = 16.5
$ws.Range("B14").Value = 6.225
$ws.Range("D15").Value = -8.231999999999999
$ws.Range("E16").Value = 16.737
$ws.Range("A18").Value = -22.049
$ws.Range("A20").Value = -20.457
$ws.Range("E20").Value = 16.598
$ws.Range("E24").Value = 17.075
$ws.Range("B26").Value = 6.097
$ws.Range("A27").Value = -21.489
$ws.Range("B27").Value = 5.62
$ws.Range("B29").Value = 6.241000000000001
$ws.Range("D33").Value = -7.955000000000001
$ws.Range("A35").Value = -19.873
$ws.Range("D35").Value = -7.551
$ws.Range("B37").Value = 8.847000000000001
$ws.Range("B38").Value = 6.198
$ws.Range("D38").Value = -8.669
$ws.Range("E39").Value = 16.222
$ws.Range("D43").Value = -7.87
$ws.Range("D44").Value = -7.741
$ws.Range("D47").Value = -7.749
$ws.Range("E47").Value = 17.01
$ws.Range("E48").Value = 17.159
$ws.Range("B51").Value = 5.326
$ws.Range("D51").Value = -8.42
$ws.Range("B52").Value = 5.11
$ws.Range("E52").Value = 16.412
$ws.Range("B55").Value = 5.705
$ws.Range("E56").Value = 16.759
$ws.Range("D57").Value = -7.972999999999999
$ws.Range("D63").Value = -7.336999999999999
$ws.Range("A69").Value = -21.404
$ws.Range("B69").Value = 6.328
$ws.Range("B70").Value = 5.412
$ws.Range("D70").Value = -6.797
$ws.Range("A76").Value = -20.66
$ws.Range("A78").Value = -19.993
$ws.Range("B81").Value = 6.02
$ws.Range("A82").Value = -22.156
$ws.Range("A83").Value = -20.219
$ws.Range("B83").Value = 7.879
$ws.Range("E84").Value = 16.617
$ws.Range("D88").Value = -7.934
$ws.Range("A93").Value = -21.609
$ws.Range("D99").Value = -8.217000000000001
$ws.Range("E100").Value = 16.489
$ws.Range("E101").Value = 16.793
$ws.Range("B102").Value = 7.13
